$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("I2").Value = 6087
$ws.Range("I3").Value = 6342
$ws.Range("I4").Value = 1464
$ws.Range("I5").Value = 591
$ws.Range("I6").Value = 7205
$ws.Range("I7").Value = 21689

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("I2").Value = 62
$ws.Range("I7").Value = 251

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("I2").Value = 69
$ws.Range("I3").Value = 55
$ws.Range("I7").Value = 235

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("I2").Value = 204
$ws.Range("I6").Value = 202
$ws.Range("I7").Value = 689

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("I2").Value = 109
$ws.Range("I3").Value = 144
$ws.Range("I7").Value = 392

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("I2").Value = 196
$ws.Range("I3").Value = 305
$ws.Range("I4").Value = 53
$ws.Range("I6").Value = 258
$ws.Range("I7").Value = 836

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("I3").Value = 53
$ws.Range("I7").Value = 216

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("I2").Value = 70
$ws.Range("I3").Value = 62
$ws.Range("I7").Value = 189

$ws = $wb.Worksheets.Item("New City")
$ws.Range("I6").Value = 147
$ws.Range("I7").Value = 508

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("I2").Value = 168
$ws.Range("I4").Value = 88
$ws.Range("I6").Value = 160
$ws.Range("I7").Value = 683
$ws.Range("I8").Value = 1302
$ws.Range("I9").Value = 107
$ws.Range("I10").Value = 153
$ws.Range("I15").Value = 248
$ws.Range("I19").Value = 596
$ws.Range("I20").Value = 547
$ws.Range("I24").Value = 62
$ws.Range("I27").Value = 190
$ws.Range("I29").Value = 1337
$ws.Range("I31").Value = 216
$ws.Range("I32").Value = 26
$ws.Range("I33").Value = 978
$ws.Range("I36").Value = 296
$ws.Range("I37").Value = 689
$ws.Range("I42").Value = 756
$ws.Range("I43").Value = 186
$ws.Range("I47").Value = 153
$ws.Range("I51").Value = 253
$ws.Range("I52").Value = 470
$ws.Range("I53").Value = 230
$ws.Range("I54").Value = 443
$ws.Range("I55").Value = 243
$ws.Range("I59").Value = 38
$ws.Range("I63").Value = 70
$ws.Range("I65").Value = 508
$ws.Range("I67").Value = 836
$ws.Range("I72").Value = 84
$ws.Range("I76").Value = 312
$ws.Range("I78").Value = 295
$ws.Range("I79").Value = 612
$ws.Range("I83").Value = 472
$ws.Range("I84").Value = 189
$ws.Range("I85").Value = 985
$ws.Range("I88").Value = 196
$ws.Range("I89").Value = 251
$ws.Range("I90").Value = 270
$ws.Range("I91").Value = 230
$ws.Range("I95").Value = 333
$ws.Range("I96").Value = 235
$ws.Range("I99").Value = 392
$ws.Range("I101").Value = 21689

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("I2").Value = 161
$ws.Range("I3").Value = 171
$ws.Range("I6").Value = 100
$ws.Range("I7").Value = 472

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("I3").Value = 121
$ws.Range("I7").Value = 333

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("I2").Value = 219
$ws.Range("I3").Value = 370
$ws.Range("I6").Value = 308
$ws.Range("I7").Value = 978

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("I2").Value = 96
$ws.Range("I5").Value = 8
$ws.Range("I6").Value = 211
$ws.Range("I7").Value = 443

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("I2").Value = 392
$ws.Range("I4").Value = 68
$ws.Range("I6").Value = 371
$ws.Range("I7").Value = 1337

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("I2").Value = 203
$ws.Range("I7").Value = 596

$ws = $wb.Worksheets.Item("River North")
$ws.Range("I2").Value = 62
$ws.Range("I7").Value = 312

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("I2").Value = 281
$ws.Range("I3").Value = 379
$ws.Range("I7").Value = 985

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("I6").Value = 46
$ws.Range("I7").Value = 160

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("I2").Value = 188
$ws.Range("I6").Value = 254
$ws.Range("I7").Value = 756

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("I6").Value = 68
$ws.Range("I7").Value = 153

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("I6").Value = 108
$ws.Range("I7").Value = 295

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("I2").Value = 73
$ws.Range("I3").Value = 77
$ws.Range("I7").Value = 243

$ws = $wb.Worksheets.Item("Dunning")
$ws.Range("I3").Value = 25
$ws.Range("I4").Value = 7
$ws.Range("I7").Value = 62

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("I2").Value = 59
$ws.Range("I3").Value = 75

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("I3").Value = 84
$ws.Range("I7").Value = 230

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("I6").Value = 182
$ws.Range("I7").Value = 612

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("I2").Value = 149
$ws.Range("I4").Value = 36
$ws.Range("I6").Value = 192
$ws.Range("I7").Value = 547

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("I3").Value = 98
$ws.Range("I7").Value = 296

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("I2").Value = 126
$ws.Range("I3").Value = 167
$ws.Range("I7").Value = 470

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("I3").Value = 44
$ws.Range("I6").Value = 52
$ws.Range("I7").Value = 153

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("I6").Value = 93
$ws.Range("I7").Value = 248

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("I2").Value = 35
$ws.Range("I7").Value = 107

$ws = $wb.Worksheets.Item("Montclare")
$ws.Range("I5").Value = 1
$ws.Range("I7").Value = 38

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("I6").Value = 35
$ws.Range("I7").Value = 168

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("I3").Value = 70
$ws.Range("I7").Value = 196

$ws = $wb.Worksheets.Item("Galewood")
$ws.Range("I3").Value = 2
$ws.Range("I7").Value = 26

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("I2").Value = 392
$ws.Range("I3").Value = 366
$ws.Range("I4").Value = 83
$ws.Range("I7").Value = 1302

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("I4").Value = 23
$ws.Range("I7").Value = 190

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("I3").Value = 68
$ws.Range("I7").Value = 270

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("I2").Value = 52
$ws.Range("I3").Value = 67
$ws.Range("I7").Value = 253

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("I4").Value = 13
$ws.Range("I6").Value = 104
$ws.Range("I7").Value = 186

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("I2").Value = 51
$ws.Range("I6").Value = 107
$ws.Range("I7").Value = 230

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("I6").Value = 40
$ws.Range("I7").Value = 84

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("I3").Value = 213
$ws.Range("I7").Value = 683

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("I6").Value = 30
$ws.Range("I7").Value = 88
